# The document has two logo pictures (Pearson/"image2.png" in the footers,
# and BTec/"image1.jpg" in the headers) whose inline-shape names were
# swapped in the original commit:
#   footers: image2.png -> image1.png
#   headers: image1.jpg -> image2.jpg
#
# InlineShape objects don't expose .Name directly until they have been
# selected, so we Select() each picture first and then rename it through
# Selection.InlineShapes(1) - this reliably reaches the shape regardless
# of which story (header/footer) it lives in.

$d = $word.ActiveDocument

function Rename-FirstInlinePicture($range, [string]$newName) {
    $shape = $range.InlineShapes.Item(1)
    $shape.Select()
    $selShape = $word.Selection.InlineShapes.Item(1)
    $selShape.Name = $newName
}

$section = $d.Sections.Item(1)

# Footers - Pearson Edexcel logo: image2.png -> image1.png
for ($i = 1; $i -le $section.Footers.Count; $i++) {
    $footer = $section.Footers.Item($i)
    if ($footer.Exists -and $footer.Range.InlineShapes.Count -gt 0) {
        Rename-FirstInlinePicture $footer.Range "image1.png"
    }
}

# Headers - BTec logo: image1.jpg -> image2.jpg
for ($i = 1; $i -le $section.Headers.Count; $i++) {
    $header = $section.Headers.Item($i)
    if ($header.Exists -and $header.Range.InlineShapes.Count -gt 0) {
        Rename-FirstInlinePicture $header.Range "image2.jpg"
    }
}
